$wb = $excel.ActiveWorkbook

# --- Typography sheet: set Wildcard Ranges (column I) for the "Large" typography (row 5) ---
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("I5").Value2 = "0-9"

# --- Translation sheet ---
# Two new rows are inserted before the old row 8 (pushing the old rows 8-12 down to 10-12),
# and three brand-new rows are appended at the end (13-15).
$ws = $wb.Worksheets.Item("Translation")

function Set-TextCell($range, $text) {
    # Excel's numeric auto-detection would turn a literal "00" into the number 0, which
    # changes both the stored type and the displayed text. Force a genuine text cell by
    # switching to a text number format while writing the value, then drop the format
    # back to General/Normal style so no visible formatting residue is left behind.
    if ($text -match '^-?\d+(\.\d+)?$') {
        $range.NumberFormat = "@"
        $range.Value2 = $text
        $range.Style = "Normal"
    } else {
        $range.Value2 = $text
    }
}

Set-TextCell $ws.Range("B4")  "SingleUseId1"
Set-TextCell $ws.Range("C4")  "Default"
Set-TextCell $ws.Range("D4")  "Center"
Set-TextCell $ws.Range("E4")  "LTR"
Set-TextCell $ws.Range("F4")  "Hour"

Set-TextCell $ws.Range("B5")  "SingleUseId2"
Set-TextCell $ws.Range("C5")  "Large"
Set-TextCell $ws.Range("D5")  "Center"
Set-TextCell $ws.Range("E5")  "LTR"
Set-TextCell $ws.Range("F5")  "<value>"

Set-TextCell $ws.Range("B6")  "SingleUseId3"
Set-TextCell $ws.Range("C6")  "Large"
Set-TextCell $ws.Range("D6")  "Left"
Set-TextCell $ws.Range("E6")  "LTR"
Set-TextCell $ws.Range("F6")  "00"

Set-TextCell $ws.Range("B7")  "SingleUseId4"
Set-TextCell $ws.Range("C7")  "Default"
Set-TextCell $ws.Range("D7")  "Center"
Set-TextCell $ws.Range("E7")  "LTR"
Set-TextCell $ws.Range("F7")  "Minute"

Set-TextCell $ws.Range("B8")  "SingleUseId5"
Set-TextCell $ws.Range("C8")  "Large"
Set-TextCell $ws.Range("D8")  "Center"
Set-TextCell $ws.Range("E8")  "LTR"
Set-TextCell $ws.Range("F8")  "<value>"

Set-TextCell $ws.Range("B9")  "SingleUseId6"
Set-TextCell $ws.Range("C9")  "Large"
Set-TextCell $ws.Range("D9")  "Left"
Set-TextCell $ws.Range("E9")  "LTR"
Set-TextCell $ws.Range("F9")  "00"

Set-TextCell $ws.Range("B10") "SingleUseId7"
Set-TextCell $ws.Range("C10") "Typography_00"
Set-TextCell $ws.Range("D10") "Center"
Set-TextCell $ws.Range("E10") "LTR"
Set-TextCell $ws.Range("F10") "Save"

Set-TextCell $ws.Range("B11") "SingleUseId8"
Set-TextCell $ws.Range("C11") "Typography_00"
Set-TextCell $ws.Range("D11") "Center"
Set-TextCell $ws.Range("E11") "LTR"
Set-TextCell $ws.Range("F11") "Save"

Set-TextCell $ws.Range("B12") "SingleUseId9"
Set-TextCell $ws.Range("C12") "Typography_00"
Set-TextCell $ws.Range("D12") "Center"
Set-TextCell $ws.Range("E12") "LTR"
Set-TextCell $ws.Range("F12") "Save"

# New rows 13-15
Set-TextCell $ws.Range("B13") "SingleUseId10"
Set-TextCell $ws.Range("C13") "Large"
Set-TextCell $ws.Range("D13") "Center"
Set-TextCell $ws.Range("E13") "LTR"
Set-TextCell $ws.Range("F13") "<hour>:<min>"

Set-TextCell $ws.Range("B14") "SingleUseId11"
Set-TextCell $ws.Range("C14") "Large"
Set-TextCell $ws.Range("D14") "Left"
Set-TextCell $ws.Range("E14") "LTR"
Set-TextCell $ws.Range("F14") "00"

Set-TextCell $ws.Range("B15") "SingleUseId12"
Set-TextCell $ws.Range("C15") "Large"
Set-TextCell $ws.Range("D15") "Left"
Set-TextCell $ws.Range("E15") "LTR"
Set-TextCell $ws.Range("F15") "00"
